$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Predictions")

function Set-TextCell {
    param($row, $col, $val)
    $c = $ws.Cells.Item($row, $col)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.NumberFormat = "General"
}

function Set-NumberCell {
    param($row, $col, $val)
    $ws.Cells.Item($row, $col).Value = $val
}

# ---- Update existing row 283: result of the match came in ----
Set-TextCell 283 12 "Completed"
Set-TextCell 283 13 "Home Win"
Set-TextCell 283 14 "Acierto"
Set-NumberCell 283 15 1.48
Set-NumberCell 283 16 45
Set-TextCell 283 17 "2025-10-10 04:27:15"

# ---- New prediction rows appended to the tracker ----
$newRows = @(
    @{ Row=284; A="2025-10-11"; B="Liga de Expansión MX"; C="Cancún";    D="Irapuato";            E="Home Win"; F="63.78%"; G=1.95; H="23.14%"; I=1.7; J=0.02566256111386165; K=0.2566256111386165 },
    @{ Row=285; A="2025-10-11"; B="Liga de Expansión MX"; C="Venados FC"; D="Tepatitlán";          E="Home Win"; F="52.44%"; G=2.2;  H="14.21%"; I=0.9; J=0.01280385768872484; K=0.1280385768872484 },
    @{ Row=286; A="2025-10-11"; B="Liga de Expansión MX"; C="Tlaxcala";    D="CDS Tampico Madero"; E="Away Win"; F="60.33%"; G=2;    H="19.45%"; I=1.4; J=0.0206530237185212;  K=0.206530237185212 },
    @{ Row=287; A="2025-10-11"; B="Major League Soccer";  C="Inter Miami"; D="Atlanta United FC";  E="Home Win"; F="87.79%"; G=1.45; H="26.02%"; I=3.3; J=0.05;                 K=0.6064470615842978 },
    @{ Row=288; A="2025-10-11"; B="Major League Soccer";  C="Orlando City SC"; D="Vancouver Whitecaps"; E="Home Win"; F="54.05%"; G=2.1;  H="12.37%"; I=0.8; J=0.0122756764106551;  K=0.122756764106551 }
)

foreach ($r in $newRows) {
    $row = $r.Row
    Set-TextCell   $row 1  $r.A
    Set-TextCell   $row 2  $r.B
    Set-TextCell   $row 3  $r.C
    Set-TextCell   $row 4  $r.D
    Set-TextCell   $row 5  $r.E
    Set-TextCell   $row 6  $r.F
    Set-NumberCell $row 7  $r.G
    Set-TextCell   $row 8  $r.H
    Set-NumberCell $row 9  $r.I
    Set-NumberCell $row 10 $r.J
    Set-NumberCell $row 11 $r.K
    Set-TextCell   $row 12 "Pending"
}
